$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# Locate the two paragraphs that need editing: "Apache Version" and the
# "PHP Version 7.1.33" paragraph that follows it.
# ----------------------------------------------------------------------
$apacheIdx = 0
$phpIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($apacheIdx -eq 0 -and $t -like "Apache Version*") {
        $apacheIdx = $i
    }
    if ($phpIdx -eq 0 -and $t -like "PHP Version*7.1.33*") {
        $phpIdx = $i
    }
}

# ----------------------------------------------------------------------
# Rewrite the "Apache Version" paragraph:
#  - drop the explicit <w:spacing .../> on the paragraph mark
#  - swap the paragraph mark's <w:highlight/> for <w:lang w:val="de-DE"/>
#  - drop the run-level highlight, keep the text but add a trailing space
#  - append the (moved) "_GoBack" bookmark followed by a new run "2.4.41"
# ----------------------------------------------------------------------
$apachePara = $d.Paragraphs.Item($apacheIdx)
$apacheXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    '<w:p w:rsidR="00717B45" w:rsidRPr="00A521D2" w:rsidRDefault="00717B45" w:rsidP="00717B45">' + `
        '<w:pPr>' + `
            '<w:pStyle w:val="Listenabsatz"/>' + `
            '<w:numPr><w:ilvl w:val="1"/><w:numId w:val="38"/></w:numPr>' + `
            '<w:rPr><w:lang w:val="de-DE"/></w:rPr>' + `
        '</w:pPr>' + `
        '<w:r w:rsidRPr="00A521D2"><w:t xml:space="preserve">Apache Version </w:t></w:r>' + `
        '<w:bookmarkStart w:id="1" w:name="_GoBack"/>' + `
        '<w:bookmarkEnd w:id="1"/>' + `
        '<w:r w:rsidRPr="00A521D2"><w:t>2.4.41</w:t></w:r>' + `
    '</w:p>' + `
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$apachePara.Range.InsertXML($apacheXml)

# Re-fetch the paragraph collection: the PHP paragraph index is unchanged
# because the Apache paragraph was replaced in place (still one paragraph).
$phpPara = $d.Paragraphs.Item($phpIdx)

# ----------------------------------------------------------------------
# Rewrite the "PHP Version 7.1.33" paragraph, removing the now-relocated
# "_GoBack" bookmark that used to sit between "Version " and "7.1.33".
# ----------------------------------------------------------------------
$phpXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    '<w:p w:rsidR="007A28BC" w:rsidRPr="007A28BC" w:rsidRDefault="007A28BC" w:rsidP="007A28BC">' + `
        '<w:pPr>' + `
            '<w:pStyle w:val="Listenabsatz"/>' + `
            '<w:numPr><w:ilvl w:val="1"/><w:numId w:val="38"/></w:numPr>' + `
            '<w:rPr><w:lang w:val="de-DE"/></w:rPr>' + `
        '</w:pPr>' + `
        '<w:r w:rsidRPr="007A28BC"><w:t xml:space="preserve">PHP </w:t></w:r>' + `
        '<w:r w:rsidR="00D63BBC"><w:t xml:space="preserve">Version </w:t></w:r>' + `
        '<w:r w:rsidRPr="007A28BC"><w:t>7.1.33</w:t></w:r>' + `
    '</w:p>' + `
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$phpPara.Range.InsertXML($phpXml)
